$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-16 10:29:58"
$wsZhCn.Range("G5").Value = "2016-02-16 10:30:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-16 10:30:13"
$wsDeDe.Range("G5").Value = "2016-02-16 10:31:22"
